$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2444.8462
$ws.Range("I32").Value = 1981.3334
$ws.Range("J32").Value = 2842.1428
$ws.Range("K32").Value = 1981.3334
$ws.Range("L32").Value = 2842.1428
$ws.Range("M32").Value = -1655.3334
$ws.Range("N32").Value = -3494.1428
$ws.Range("H33").Value = 5716.9473
$ws.Range("I33").Value = 6507.8125
$ws.Range("K33").Value = 6507.8125
$ws.Range("M33").Value = -6278.8125
$ws.Range("H51").Value = 7677
$ws.Range("I51").Value = 8475.75
$ws.Range("J51").Value = 6399
$ws.Range("K51").Value = 8475.75
$ws.Range("L51").Value = 6399
$ws.Range("M51").Value = -7991.75
$ws.Range("N51").Value = -7367
$ws.Range("H80").Value = 498.7857
$ws.Range("I80").Value = 271.66666
$ws.Range("K80").Value = 814.9999799999999
$ws.Range("M80").Value = 183.0000200000001
$ws.Range("H83").Value = 498.7857
$ws.Range("I83").Value = 271.66666
$ws.Range("K83").Value = 2444.99994
$ws.Range("M83").Value = 2547.00006
$ws.Range("H100").Value = 3089.5881
$ws.Range("I100").Value = 1653.3334
$ws.Range("J100").Value = 3873
$ws.Range("K100").Value = 1653.3334
$ws.Range("L100").Value = 3873
$ws.Range("M100").Value = -1112.3334
$ws.Range("N100").Value = -4955
$ws.Range("H101").Value = 599.375
$ws.Range("I101").Value = 499.33334
$ws.Range("K101").Value = 1498.00002
$ws.Range("M101").Value = 123.9999800000001
$ws.Range("H106").Value = 8571.3125
$ws.Range("I106").Value = 7040.4287
$ws.Range("J106").Value = 19287.5
$ws.Range("K106").Value = 7040.4287
$ws.Range("L106").Value = 19287.5
$ws.Range("M106").Value = -6409.4287
$ws.Range("N106").Value = -20549.5
$ws.Range("H129").Value = 1293.25
$ws.Range("J129").Value = 1249.625
$ws.Range("L129").Value = 3748.875
$ws.Range("N129").Value = -13748.875
$ws.Range("H141").Value = 5187
$ws.Range("I141").Value = 5187
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 15561
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -10381
$ws.Range("N141").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13976.326
$ws.Range("I32").Value = 14011.591
$ws.Range("K32").Value = 14011.591
$ws.Range("M32").Value = -13724.591
$ws.Range("H61").Value = 1576.3334
$ws.Range("I61").Value = 1576.3334
$ws.Range("K61").Value = 1576.3334
$ws.Range("M61").Value = -1364.3334
$ws.Range("H74").Value = 22276.883
$ws.Range("I74").Value = 23749.639
$ws.Range("J74").Value = 4972
$ws.Range("K74").Value = 23749.639
$ws.Range("L74").Value = 4972
$ws.Range("M74").Value = -22875.639
$ws.Range("N74").Value = -6720
$ws.Range("H77").Value = 22276.883
$ws.Range("I77").Value = 23749.639
$ws.Range("J77").Value = 4972
$ws.Range("K77").Value = 118748.195
$ws.Range("L77").Value = 24860
$ws.Range("M77").Value = -114380.195
$ws.Range("N77").Value = -33596
$ws.Range("H132").Value = 33625.156
$ws.Range("I132").Value = 35733.535
$ws.Range("J132").Value = 1999.5
$ws.Range("K132").Value = 107200.605
$ws.Range("L132").Value = 5998.5
$ws.Range("M132").Value = -104670.605
$ws.Range("N132").Value = -11058.5
$ws.Range("H136").Value = 1576.3334
$ws.Range("I136").Value = 1576.3334
$ws.Range("K136").Value = 4729.0002
$ws.Range("M136").Value = -2179.0002
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3422.65
$ws.Range("I86").Value = 3144.353
$ws.Range("K86").Value = 3144.353
$ws.Range("M86").Value = -2021.353
$ws.Range("H89").Value = 3422.65
$ws.Range("I89").Value = 3144.353
$ws.Range("K89").Value = 15721.765
$ws.Range("M89").Value = -10105.765
$ws.Range("H134").Value = 2457.2
$ws.Range("I134").Value = 2348.8276
$ws.Range("K134").Value = 7046.4828
$ws.Range("M134").Value = -4511.4828
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1392.2307
$ws.Range("I22").Value = 337.42856
$ws.Range("K22").Value = 337.42856
$ws.Range("M22").Value = 12.57144
$ws.Range("H31").Value = 2763.0454
$ws.Range("I31").Value = 2643.7646
$ws.Range("J31").Value = 3168.6
$ws.Range("K31").Value = 2643.7646
$ws.Range("L31").Value = 3168.6
$ws.Range("M31").Value = -2348.7646
$ws.Range("N31").Value = -3758.6
$ws.Range("H34").Value = 2763.0454
$ws.Range("I34").Value = 2643.7646
$ws.Range("J34").Value = 3168.6
$ws.Range("K34").Value = 2643.7646
$ws.Range("L34").Value = 3168.6
$ws.Range("M34").Value = -2441.7646
$ws.Range("N34").Value = -3572.6
$ws.Range("H58").Value = 64127.562
$ws.Range("I58").Value = 72795.78999999999
$ws.Range("K58").Value = 72795.78999999999
$ws.Range("M58").Value = -72592.78999999999
$ws.Range("H105").Value = 698
$ws.Range("I105").Value = 698
$ws.Range("K105").Value = 698
$ws.Range("M105").Value = 1049
$ws.Range("H132").Value = 2683.8235
$ws.Range("I132").Value = 2570.4375
$ws.Range("K132").Value = 7711.3125
$ws.Range("M132").Value = -5181.3125
$ws.Range("H136").Value = 64127.562
$ws.Range("I136").Value = 72795.78999999999
$ws.Range("K136").Value = 218387.37
$ws.Range("M136").Value = -215837.37
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 492.5
$ws.Range("J12").Value = 288.1111
$ws.Range("L12").Value = 864.3333
$ws.Range("N12").Value = -1210.3333
$ws.Range("H68").Value = 999
$ws.Range("I68").Value = 999
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2997
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("M68").Value = -2186
$ws.Range("H71").Value = 999
$ws.Range("I71").Value = 999
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 8991
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("M71").Value = -4935
$ws.Range("H97").Value = 380.5
$ws.Range("I97").Value = 173.2
$ws.Range("K97").Value = 519.5999999999999
$ws.Range("M97").Value = -23.59999999999991
$ws.Range("H137").Value = 4154.364
$ws.Range("J137").Value = 5124.875
$ws.Range("L137").Value = 15374.625
$ws.Range("N137").Value = -25574.625
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 49083
$ws.Range("J15").Value = 49083
$ws.Range("L15").Value = 49083
$ws.Range("N15").Value = -49659
$ws.Range("H81").Value = 49083
$ws.Range("J81").Value = 49083
$ws.Range("L81").Value = 49083
$ws.Range("N81").Value = -51079
$ws.Range("H84").Value = 49083
$ws.Range("J84").Value = 49083
$ws.Range("L84").Value = 147249
$ws.Range("N84").Value = -157233
$ws.Range("H97").Value = 751.3333
$ws.Range("I97").Value = 579.1818
$ws.Range("J97").Value = 1224.75
$ws.Range("K97").Value = 579.1818
$ws.Range("L97").Value = 1224.75
$ws.Range("M97").Value = -83.18179999999995
$ws.Range("N97").Value = -2216.75
$ws.Range("H126").Value = 8000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 24000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -28940
$ws.Range("H132").Value = 50691.617
$ws.Range("I132").Value = 65514.75
$ws.Range("J132").Value = 3257.6
$ws.Range("K132").Value = 196544.25
$ws.Range("L132").Value = 9772.799999999999
$ws.Range("M132").Value = -194014.25
$ws.Range("N132").Value = -14832.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1739.7
$ws.Range("I22").Value = 1388.6666
$ws.Range("K22").Value = 1388.6666
$ws.Range("M22").Value = -1093.6666
$ws.Range("H27").Value = 1739.7
$ws.Range("I27").Value = 1388.6666
$ws.Range("K27").Value = 1388.6666
$ws.Range("M27").Value = -1281.6666
$ws.Range("H93").Value = 2012.9286
$ws.Range("I93").Value = 1683.6666
$ws.Range("J93").Value = 2168.8948
$ws.Range("K93").Value = 1683.6666
$ws.Range("L93").Value = 2168.8948
$ws.Range("M93").Value = -435.6666
$ws.Range("N93").Value = -4664.8948
$ws.Range("H100").Value = 4324.9
$ws.Range("I100").Value = 3562.5
$ws.Range("K100").Value = 3562.5
$ws.Range("M100").Value = -3021.5
$ws.Range("H132").Value = 55271.652
$ws.Range("I132").Value = 65776.52
$ws.Range("J132").Value = 5373.5
$ws.Range("K132").Value = 197329.56
$ws.Range("L132").Value = 16120.5
$ws.Range("M132").Value = -194799.56
$ws.Range("N132").Value = -21180.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 45849.5
$ws.Range("I45").Value = 7999.5
$ws.Range("J45").Value = 64774.5
$ws.Range("K45").Value = 7999.5
$ws.Range("L45").Value = 64774.5
$ws.Range("M45").Value = -7508.5
$ws.Range("N45").Value = -65756.5
$ws.Range("H81").Value = 3787.25
$ws.Range("J81").Value = 15333.667
$ws.Range("L81").Value = 30667.334
$ws.Range("N81").Value = -32789.334
$ws.Range("H84").Value = 3787.25
$ws.Range("J84").Value = 15333.667
$ws.Range("L84").Value = 153336.67
$ws.Range("N84").Value = -163944.67
$ws.Range("H132").Value = 52537.27
$ws.Range("I132").Value = 57272.46
$ws.Range("K132").Value = 171817.38
$ws.Range("M132").Value = -169287.38
